$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.511.45"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "2.489.49"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "2.487.68"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "2.944.89"
$ws.Range("E14").Value = "  -0.81%  "
$ws.Range("D15").Value = "69.262.61"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("D18").Value = "2.493.91"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "347.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("E23").Value = "  -4.27%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.57%  "
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "2.616.08"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  -4.32%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("E32").Value = "  -4.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "434.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("E35").Value = "  -1.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "155.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("E42").Value = "  -3.56%  "
$ws.Range("B43").Value = "POPCAT"
$ws.Range("C43").Value = "https://coinranking.com/coin/sLBuDEsp6+popcat-popcat"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +51.63%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("E45").Value = "  -4.41%  "
$ws.Range("E46").Value = "  -5.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "138.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("E48").Value = "  -1.59%  "
$ws.Range("E49").Value = "  -4.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0722"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("E51").Value = "  -0.89%  "
